$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2999
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null

$ws.Range("H76").Value = 3150
$ws.Range("I76").Value = 3051.8572
$ws.Range("K76").Value = 3051.8572
$ws.Range("M76").Value = -2736.8572

$ws.Range("H79").Value = 3150
$ws.Range("I79").Value = 3051.8572
$ws.Range("K79").Value = 3051.8572
$ws.Range("M79").Value = -1959.8572

$ws.Range("H100").Value = 2352
$ws.Range("I100").Value = 1469.3334
$ws.Range("K100").Value = 1469.3334
$ws.Range("M100").Value = -928.3334

$ws.Range("H106").Value = 2832.6667
$ws.Range("I106").Value = 2832.6667
$ws.Range("K106").Value = 2832.6667
$ws.Range("M106").Value = -2201.6667

$ws.Range("H108").Value = 99697.78
$ws.Range("J108").Value = 99697.78
$ws.Range("L108").Value = 99697.78
$ws.Range("N108").Value = -107377.78

$ws.Range("H110").Value = 67859
$ws.Range("J110").Value = 67859
$ws.Range("L110").Value = 67859
$ws.Range("N110").Value = -76039

$ws.Range("H120").Value = 50248
$ws.Range("J120").Value = 50248
$ws.Range("L120").Value = 50248
$ws.Range("N120").Value = -59924

$ws.Range("H123").Value = 84168.57000000001
$ws.Range("J123").Value = 85696.664
$ws.Range("L123").Value = 85696.664
$ws.Range("N123").Value = -95496.664

$ws.Range("H125").Value = 4083.25
$ws.Range("I125").Value = 3666.5
$ws.Range("J125").Value = 4500
$ws.Range("K125").Value = 32998.5
$ws.Range("L125").Value = 40500
$ws.Range("M125").Value = -30538.5
$ws.Range("N125").Value = -45420

$ws.Range("H133").Value = 74975.55499999999
$ws.Range("J133").Value = 74975.55499999999
$ws.Range("L133").Value = 74975.55499999999
$ws.Range("N133").Value = -85095.55499999999

$ws.Range("H136").Value = 77977.14
$ws.Range("J136").Value = 77977.14
$ws.Range("L136").Value = 77977.14
$ws.Range("N136").Value = -88177.14

$ws.Range("H139").Value = 71623.86
$ws.Range("J139").Value = 71623.86
$ws.Range("L139").Value = 71623.86
$ws.Range("N139").Value = -81903.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9619.227999999999
$ws.Range("I32").Value = 4268.75
$ws.Range("K32").Value = 4268.75
$ws.Range("M32").Value = -3981.75

$ws.Range("H74").Value = 35680
$ws.Range("I74").Value = 54156.895
$ws.Range("K74").Value = 54156.895
$ws.Range("M74").Value = -53282.895

$ws.Range("H77").Value = 35680
$ws.Range("I77").Value = 54156.895
$ws.Range("K77").Value = 270784.475
$ws.Range("M77").Value = -266416.475

$ws.Range("H121").Value = 82202.125
$ws.Range("J121").Value = 82202.125
$ws.Range("L121").Value = 82202.125
$ws.Range("N121").Value = -85696.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2512.8823
$ws.Range("I94").Value = 1623.1428
$ws.Range("K94").Value = 1623.1428
$ws.Range("M94").Value = -1172.1428

$ws.Range("H99").Value = 3584418.2
$ws.Range("I99").Value = 143680.72
$ws.Range("J99").Value = 15627000
$ws.Range("K99").Value = 143680.72
$ws.Range("L99").Value = 15627000
$ws.Range("M99").Value = -142182.72
$ws.Range("N99").Value = -15629996

$ws.Range("H110").Value = 51670.9
$ws.Range("J110").Value = 51670.9
$ws.Range("L110").Value = 51670.9
$ws.Range("N110").Value = -59850.9

$ws.Range("H114").Value = 99970.336
$ws.Range("J114").Value = 99970.336
$ws.Range("L114").Value = 99970.336
$ws.Range("N114").Value = -108648.336

$ws.Range("H117").Value = 99969.42999999999
$ws.Range("J117").Value = 99969.42999999999
$ws.Range("L117").Value = 99969.42999999999
$ws.Range("N117").Value = -109147.43

$ws.Range("H118").Value = 71579.336

$ws.Range("H122").Value = 62844.375
$ws.Range("J122").Value = 62844.375
$ws.Range("L122").Value = 62844.375
$ws.Range("N122").Value = -72644.375

$ws.Range("H132").Value = 32384.154
$ws.Range("J132").Value = 32384.154
$ws.Range("L132").Value = 32384.154
$ws.Range("N132").Value = -42504.15399999999

$ws.Range("H135").Value = 21000
$ws.Range("J135").Value = 21000
$ws.Range("L135").Value = 21000
$ws.Range("N135").Value = -31140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1000000
$ws.Range("J6").Value = 1000000
$ws.Range("L6").Value = 1000000
$ws.Range("N6").Value = -1000226

$ws.Range("H108").Value = 40744.535
$ws.Range("J108").Value = 40744.535
$ws.Range("L108").Value = 40744.535
$ws.Range("N108").Value = -48424.535

$ws.Range("H116").Value = 72203.2
$ws.Range("J116").Value = 72203.2
$ws.Range("L116").Value = 72203.2
$ws.Range("N116").Value = -81381.2

$ws.Range("H134").Value = 3020096.2
$ws.Range("I134").Value = 3249048.8
$ws.Range("K134").Value = 9747146.399999999
$ws.Range("M134").Value = -9744611.399999999

$ws.Range("H141").Value = 93325
$ws.Range("I141").Value = 60000
$ws.Range("K141").Value = 60000
$ws.Range("M141").Value = -54820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 185
$ws.Range("I24").Value = 180
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 540
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = -310
$ws.Range("N24").Value = -1060

$ws.Range("H37").Value = 67776
$ws.Range("J37").Value = 67776
$ws.Range("L37").Value = 203328
$ws.Range("N37").Value = -203552

$ws.Range("H57").Value = 5333
$ws.Range("I57").Value = 3999
$ws.Range("J57").Value = 6000
$ws.Range("K57").Value = 11997
$ws.Range("L57").Value = 18000
$ws.Range("M57").Value = -11438
$ws.Range("N57").Value = -19118

$ws.Range("H131").Value = 126443.625
$ws.Range("I131").Value = 167803.17
$ws.Range("K131").Value = 503409.51
$ws.Range("M131").Value = -498369.51

$ws.Range("H138").Value = 4948.25
$ws.Range("I138").Value = 3396.5
$ws.Range("K138").Value = 10189.5
$ws.Range("M138").Value = -5049.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 764.5263
$ws.Range("J2").Value = 241.75
$ws.Range("L2").Value = 241.75
$ws.Range("N2").Value = -467.75

$ws.Range("H20").Value = 50138.5
$ws.Range("J20").Value = 50138.5
$ws.Range("L20").Value = 50138.5
$ws.Range("N20").Value = -50628.5

$ws.Range("H24").Value = 34353.816
$ws.Range("I24").Value = 30006
$ws.Range("J24").Value = 34788.6
$ws.Range("K24").Value = 30006
$ws.Range("L24").Value = 34788.6
$ws.Range("M24").Value = -29833
$ws.Range("N24").Value = -35134.6

$ws.Range("H70").Value = 5391.923
$ws.Range("I70").Value = 5179.5
$ws.Range("J70").Value = 6100
$ws.Range("K70").Value = 5179.5
$ws.Range("L70").Value = 6100
$ws.Range("M70").Value = -4909.5
$ws.Range("N70").Value = -6640

$ws.Range("H73").Value = 5391.923
$ws.Range("I73").Value = 5179.5
$ws.Range("J73").Value = 6100
$ws.Range("K73").Value = 5179.5
$ws.Range("L73").Value = 6100
$ws.Range("M73").Value = -4243.5
$ws.Range("N73").Value = -7972

$ws.Range("H93").Value = 17024.334
$ws.Range("J93").Value = 17024.334
$ws.Range("L93").Value = 17024.334
$ws.Range("N93").Value = -20768.334

$ws.Range("H109").Value = 29028.5
$ws.Range("I109").Value = 8250
$ws.Range("K109").Value = 8250
$ws.Range("M109").Value = -7210

$ws.Range("H114").Value = 73626.44500000001
$ws.Range("J114").Value = 73626.44500000001
$ws.Range("L114").Value = 73626.44500000001
$ws.Range("N114").Value = -82304.44500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 703534.7
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null

$ws.Range("H71").Value = 703534.7
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null

$ws.Range("H136").Value = 3643.0715
$ws.Range("I136").Value = 3802.35
$ws.Range("K136").Value = 11407.05
$ws.Range("M136").Value = -8857.049999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 16000
$ws.Range("J105").Value = 16000
$ws.Range("L105").Value = 16000
$ws.Range("N105").Value = -22988

$ws.Range("H121").Value = 99995
$ws.Range("J121").Value = 99995
$ws.Range("L121").Value = 99995
$ws.Range("N121").Value = -103489

$ws.Range("H122").Value = 2139.3333
$ws.Range("I122").Value = 1344.1428
$ws.Range("J122").Value = 2835.125
$ws.Range("K122").Value = 4032.4284
$ws.Range("L122").Value = 8505.375
$ws.Range("M122").Value = -1582.4284
$ws.Range("N122").Value = -13405.375

$ws.Range("H126").Value = 280751
$ws.Range("I126").Value = 553502
$ws.Range("K126").Value = 1660506
$ws.Range("M126").Value = -1658036
